$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 11; $r -le 17; $r++) {
    $ws.Cells.Item($r, 7).Value = "55,110,0"
}

$ws.Range("G15").Select()
